# Netzplan Übung 3 FIAE D ++
# - Fill in the "D" (Dauer/duration) column on Tabelle2 (D3:D8)
# - Move the active sheet / selection from Tabelle1!L4 to Tabelle2,
#   leaving a new selection on Tabelle1 at AD19 and selecting F8 on Tabelle2

$wb = $excel.ActiveWorkbook

# Tabelle1 keeps its selection moved to AD19, but is no longer the active tab.
$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws1.Activate()
$ws1.Range("AD19").Select()

# Tabelle2 becomes the active sheet; fill in the duration values and select F8.
$ws2 = $wb.Worksheets.Item("Tabelle2")
$ws2.Activate()

$ws2.Range("D3").Value = 4
$ws2.Range("D4").Value = 3
$ws2.Range("D5").Value = 2
$ws2.Range("D6").Value = 9
$ws2.Range("D7").Value = 2
$ws2.Range("D8").Value = 6

$ws2.Range("F8").Select()
